$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 474.69232
$ws.Range("I33").Value = 502
$ws.Range("J33").Value = 324.5
$ws.Range("K33").Value = 502
$ws.Range("L33").Value = 324.5
$ws.Range("M33").Value = -273
$ws.Range("N33").Value = -782.5

$ws.Range("H40").Value = 3463.5386
$ws.Range("I40").Value = 3053
$ws.Range("J40").Value = 4387.25
$ws.Range("K40").Value = 3053
$ws.Range("L40").Value = 4387.25
$ws.Range("M40").Value = -2878
$ws.Range("N40").Value = -4737.25

$ws.Range("H70").Value = 4364.6665
$ws.Range("I70").Value = 6916.6665
$ws.Range("J70").Value = 2663.3333
$ws.Range("K70").Value = 20749.9995
$ws.Range("L70").Value = 7989.999899999999
$ws.Range("M70").Value = -20479.9995
$ws.Range("N70").Value = -8529.999899999999

$ws.Range("H73").Value = 4364.6665
$ws.Range("I73").Value = 6916.6665
$ws.Range("J73").Value = 2663.3333
$ws.Range("K73").Value = 20749.9995
$ws.Range("L73").Value = 7989.999899999999
$ws.Range("M73").Value = -19813.9995
$ws.Range("N73").Value = -9861.999899999999

$ws.Range("H92").Value = 53363.42
$ws.Range("I92").Value = 77492.766
$ws.Range("K92").Value = 77492.766
$ws.Range("M92").Value = -76244.766

$ws.Range("H106").Value = 7426.9165
$ws.Range("I106").Value = 6829.4546
$ws.Range("K106").Value = 6829.4546
$ws.Range("M106").Value = -6198.4546

$ws.Range("H112").Value = 4833.3335
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 4833.3335
$ws.Range("K112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("M112").Value = 14500.0005
$ws.Range("N112").Value = -16716.0005

$ws.Range("H135").Value = 555.8
$ws.Range("I135").Value = 562.2917
$ws.Range("K135").Value = 5060.6253
$ws.Range("M135").Value = -2525.6253

$ws.Range("H137").Value = 2040.2858
$ws.Range("I137").Value = 2040.2858
$ws.Range("K137").Value = 6120.857400000001
$ws.Range("M137").Value = -3570.857400000001

$ws.Range("H141").Value = 1464.871
$ws.Range("I141").Value = 1464.871
$ws.Range("K141").Value = 4394.613
$ws.Range("M141").Value = 785.3869999999997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 17626514
$ws.Range("I11").Value = 17626514
$ws.Range("K11").Value = 17626514
$ws.Range("M11").Value = -17626370

$ws.Range("H32").Value = 14756.833
$ws.Range("I32").Value = 14756.833
$ws.Range("K32").Value = 14756.833
$ws.Range("M32").Value = -14469.833

$ws.Range("H61").Value = 1753.7273
$ws.Range("I61").Value = 929.1
$ws.Range("J61").Value = 10000
$ws.Range("K61").Value = 929.1
$ws.Range("L61").Value = 10000
$ws.Range("M61").Value = -717.1
$ws.Range("N61").Value = -10424

$ws.Range("H74").Value = 2397.3333
$ws.Range("I74").Value = 2397.3333
$ws.Range("K74").Value = 2397.3333
$ws.Range("M74").Value = -1523.3333

$ws.Range("H77").Value = 2397.3333
$ws.Range("I77").Value = 2397.3333
$ws.Range("K77").Value = 11986.6665
$ws.Range("M77").Value = -7618.666499999999

$ws.Range("H97").Value = 1611.2667
$ws.Range("I97").Value = 1051.4615
$ws.Range("K97").Value = 1051.4615
$ws.Range("M97").Value = -555.4614999999999

$ws.Range("H110").Value = 4046.652
$ws.Range("I110").Value = 3916.3
$ws.Range("J110").Value = 4915.6665
$ws.Range("K110").Value = 3916.3
$ws.Range("L110").Value = 4915.6665
$ws.Range("M110").Value = -1871.3
$ws.Range("N110").Value = -9005.666499999999

$ws.Range("H132").Value = 1538.3846
$ws.Range("I132").Value = 1399.92
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 4199.76
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -1669.76
$ws.Range("N132").Value = -20060

$ws.Range("H136").Value = 1753.7273
$ws.Range("I136").Value = 929.1
$ws.Range("J136").Value = 10000
$ws.Range("K136").Value = 2787.3
$ws.Range("L136").Value = 30000
$ws.Range("M136").Value = -237.3000000000002
$ws.Range("N136").Value = -35100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2022.9166
$ws.Range("I20").Value = 1808.6666
$ws.Range("K20").Value = 1808.6666
$ws.Range("M20").Value = -1561.6666

$ws.Range("H86").Value = 6504.1934
$ws.Range("I86").Value = 4029.682
$ws.Range("K86").Value = 4029.682
$ws.Range("M86").Value = -2906.682

$ws.Range("H89").Value = 6504.1934
$ws.Range("I89").Value = 4029.682
$ws.Range("K89").Value = 20148.41
$ws.Range("M89").Value = -14532.41

$ws.Range("H94").Value = 597.2857
$ws.Range("I94").Value = 627.4761999999999
$ws.Range("J94").Value = 506.7143
$ws.Range("K94").Value = 627.4761999999999
$ws.Range("L94").Value = 506.7143
$ws.Range("M94").Value = -176.4761999999999
$ws.Range("N94").Value = -1408.7143

$ws.Range("H105").Value = 2833.3333
$ws.Range("I105").Value = 1750
$ws.Range("K105").Value = 1750
$ws.Range("M105").Value = -3

$ws.Range("H134").Value = 6743.926
$ws.Range("I134").Value = 6632.3335
$ws.Range("K134").Value = 19897.0005
$ws.Range("M134").Value = -17362.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 714906.1
$ws.Range("I19").Value = 909344.2
$ws.Range("J19").Value = 1966.6666
$ws.Range("K19").Value = 909344.2
$ws.Range("L19").Value = 1966.6666
$ws.Range("M19").Value = -909174.2
$ws.Range("N19").Value = -2306.6666

$ws.Range("H24").Value = 714906.1
$ws.Range("I24").Value = 909344.2
$ws.Range("J24").Value = 1966.6666
$ws.Range("K24").Value = 909344.2
$ws.Range("L24").Value = 1966.6666
$ws.Range("M24").Value = -909174.2
$ws.Range("N24").Value = -2306.6666

$ws.Range("H31").Value = 4999.5
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 4999.5
$ws.Range("K31").Value = 0
$ws.Range("L31").ClearContents()
$ws.Range("M31").Value = 4999.5
$ws.Range("N31").Value = -5589.5

$ws.Range("H34").Value = 4999.5
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 4999.5
$ws.Range("K34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("M34").Value = 4999.5
$ws.Range("N34").Value = -5403.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3478.111
$ws.Range("I5").Value = 3600.5
$ws.Range("K5").Value = 10801.5
$ws.Range("M5").Value = -10689.5

$ws.Range("H6").Value = 353.5
$ws.Range("I6").Value = 353.5
$ws.Range("K6").Value = 1060.5
$ws.Range("M6").Value = -947.5

$ws.Range("H97").Value = 1590.4
$ws.Range("I97").Value = 1599.25
$ws.Range("J97").Value = 1555
$ws.Range("K97").Value = 4797.75
$ws.Range("L97").Value = 4665
$ws.Range("M97").Value = -4301.75
$ws.Range("N97").Value = -5657

$ws.Range("H135").Value = 3478.111
$ws.Range("I135").Value = 3600.5
$ws.Range("K135").Value = 32404.5
$ws.Range("M135").Value = -29869.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9449.5
$ws.Range("I70").Value = 9339.333000000001
$ws.Range("K70").Value = 9339.333000000001
$ws.Range("M70").Value = -9069.333000000001

$ws.Range("H73").Value = 9449.5
$ws.Range("I73").Value = 9339.333000000001
$ws.Range("K73").Value = 9339.333000000001
$ws.Range("M73").Value = -8403.333000000001

$ws.Range("H97").Value = 835.9091
$ws.Range("I97").Value = 819.5
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 819.5
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -323.5
$ws.Range("N97").Value = -1992

$ws.Range("H126").Value = 3513.5557
$ws.Range("I126").Value = 3577.875
$ws.Range("K126").Value = 10733.625
$ws.Range("M126").Value = -8263.625

$ws.Range("H132").Value = 3304.35
$ws.Range("I132").Value = 2467.5386
$ws.Range("K132").Value = 7402.6158
$ws.Range("M132").Value = -4872.6158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2601.4375
$ws.Range("I68").Value = 2533.6365
$ws.Range("K68").Value = 2533.6365
$ws.Range("M68").Value = -1784.6365

$ws.Range("H71").Value = 2601.4375
$ws.Range("I71").Value = 2533.6365
$ws.Range("K71").Value = 12668.1825
$ws.Range("M71").Value = -8924.182500000001

$ws.Range("H82").Value = 1238.85
$ws.Range("I82").Value = 844.4545000000001
$ws.Range("J82").Value = 1720.8889
$ws.Range("K82").Value = 844.4545000000001
$ws.Range("L82").Value = 1720.8889
$ws.Range("M82").Value = -483.4545000000001
$ws.Range("N82").Value = -2442.8889

$ws.Range("H85").Value = 1238.85
$ws.Range("I85").Value = 844.4545000000001
$ws.Range("J85").Value = 1720.8889
$ws.Range("K85").Value = 844.4545000000001
$ws.Range("L85").Value = 1720.8889
$ws.Range("M85").Value = 403.5454999999999
$ws.Range("N85").Value = -4216.8889

$ws.Range("H122").Value = 2289.2104
$ws.Range("I122").Value = 2076.625
$ws.Range("K122").Value = 6229.875
$ws.Range("M122").Value = -3779.875
